$wb = $excel.ActiveWorkbook

# --- Workbook-level path metadata (cosmetic, matches author's machine) ---
# Not exposed via the Excel object model in a meaningful way; skip.

# --- Sheet 1: "Use case template" ---
$ws1 = $wb.Worksheets.Item("Use case template")

# Rename header labels: Baseline-SDMFFP1 -> SDMFFP1, Baseline-SDMFFP2 -> SDMFFP2
$ws1.Range("B1").Value = "SDMFFP1"
$ws1.Range("C1").Value = "SDMFFP2"

# Clear the old D column values, then delete column D (was Baseline-SDMFFP3)
$ws1.Columns.Item(4).Delete()

# Mark the new assessment results (add_assessment_result_to_matrix)
$ws1.Range("D2").Value = $null
$ws1.Range("B7").Value = 1
$ws1.Range("B17").Value = 1
$ws1.Range("B21").Value = 1
$ws1.Range("C22").Value = 1
$ws1.Range("C24").Value = 1

$ws1.Range("D24").Select()
